$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two existing strings in row 7 (D7/E7)
$ws.Range("D7").Value = "Phone num (new)Australian "
$ws.Range("E7").Value = "It should  be set as +61 Australian code up 8 to 20 characters"

# Fill in previously empty Actual output / Result cells for row 7
$ws.Range("F7").Value = "It should gets displayed "
$ws.Range("G7").Value = "Pass"

# Update the view state: scrolled/selected position
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F8").Select()
